$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")
$ws.Activate()

$ws.Range("D10").Value = "Medium"
$ws.Range("E10").Value = "dp"
$ws.Range("F10").Value = "Medium"
$ws.Range("G10").Value = "recursion"
$ws.Range("H10").Value = 45

$ws.Range("G20").Select()
